$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Sending cluster) for rows 2-7 changes from "Inflammatory-Mac" to "Neutrophils"
# (rows 8-13 remain "Resolving-Mac", so no change needed there)
$ws.Range("A2").Value = "Neutrophils"
$ws.Range("A3").Value = "Neutrophils"
$ws.Range("A4").Value = "Neutrophils"
$ws.Range("A5").Value = "Neutrophils"
$ws.Range("A6").Value = "Neutrophils"
$ws.Range("A7").Value = "Neutrophils"

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3674653333333333
$ws.Range("H2").Value = 1.102396
$ws.Range("I2").Value = 0.5572612813719676
$ws.Range("J2").Value = 0.5572612813719677
$ws.Range("M2").Value = 240.859306
$ws.Range("N2").Value = 481.718612
$ws.Range("O2").Value = 0.770495783518506
$ws.Range("P2").Value = 0.7103085338816119
$ws.Range("Q2").Value = 88.50744516572533
$ws.Range("R2").Value = 531.044670994352
$ws.Range("S2").Value = 0.4293674676152208
$ws.Range("T2").Value = 0.3958274437603108

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3674653333333333
$ws.Range("H3").Value = 1.102396
$ws.Range("I3").Value = 0.5572612813719676
$ws.Range("J3").Value = 0.5572612813719677
$ws.Range("O3").Value = 0.133770027207319
$ws.Range("P3").Value = 0.1849808797181794
$ws.Range("Q3").Value = 15.36626624198133
$ws.Range("R3").Value = 138.296396177832
$ws.Range("S3").Value = 0.07454485677071358
$ws.Range("T3").Value = 0.1030826820610665

# Row 4
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3674653333333333
$ws.Range("H4").Value = 1.102396
$ws.Range("I4").Value = 0.5572612813719676
$ws.Range("J4").Value = 0.5572612813719677
$ws.Range("M4").Value = 3.307267666666667
$ws.Range("N4").Value = 9.921803
$ws.Range("O4").Value = 0.01057976888853842
$ws.Range("P4").Value = 0.01462999594956937
$ws.Range("Q4").Value = 1.215306215554222
$ws.Range("R4").Value = 10.937755939988
$ws.Range("S4").Value = 0.005895695567446195
$ws.Range("T4").Value = 0.008152730289323725

# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.3674653333333333
$ws.Range("H5").Value = 1.102396
$ws.Range("I5").Value = 0.5572612813719676
$ws.Range("J5").Value = 0.5572612813719677
$ws.Range("M5").Value = 18.767532
$ws.Range("N5").Value = 37.535064
$ws.Range("O5").Value = 0.06003631129389966
$ws.Range("P5").Value = 0.05534657705729765
$ws.Range("Q5").Value = 6.896417402223999
$ws.Range("R5").Value = 41.378504413344
$ws.Range("S5").Value = 0.03345591176048485
$ws.Range("T5").Value = 0.03084250445050204

# Row 6
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.3674653333333333
$ws.Range("H6").Value = 1.102396
$ws.Range("I6").Value = 0.5572612813719676
$ws.Range("J6").Value = 0.5572612813719677
$ws.Range("M6").Value = 1.967337333333333
$ws.Range("N6").Value = 5.902012
$ws.Range("O6").Value = 0.006293404831498911
$ws.Range("P6").Value = 0.008702693618721296
$ws.Range("Q6").Value = 0.7229282689724443
$ws.Range("R6").Value = 6.506354420751999
$ws.Range("S6").Value = 0.003507070840593615
$ws.Range("T6").Value = 0.004849674197356276

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.3674653333333333
$ws.Range("H7").Value = 1.102396
$ws.Range("I7").Value = 0.5572612813719676
$ws.Range("J7").Value = 0.5572612813719677
$ws.Range("M7").Value = 5.884659333333333
$ws.Range("N7").Value = 17.653978
$ws.Range("O7").Value = 0.01882470426023795
$ws.Range("P7").Value = 0.02603131977462027
$ws.Range("Q7").Value = 2.162408303476444
$ws.Range("R7").Value = 19.461674731288
$ws.Range("S7").Value = 0.01049027881750854
$ws.Range("T7").Value = 0.01450624661340833

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.2919476666666667
$ws.Range("H8").Value = 0.875843
$ws.Range("I8").Value = 0.4427387186280323
$ws.Range("J8").Value = 0.4427387186280324
$ws.Range("M8").Value = 240.859306
$ws.Range("N8").Value = 481.718612
$ws.Range("O8").Value = 0.770495783518506
$ws.Range("P8").Value = 0.7103085338816119
$ws.Range("Q8").Value = 70.31831238165266
$ws.Range("R8").Value = 421.909874289916
$ws.Range("S8").Value = 0.3411283159032851
$ws.Range("T8").Value = 0.3144810901213012

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.2919476666666667
$ws.Range("H9").Value = 0.875843
$ws.Range("I9").Value = 0.4427387186280323
$ws.Range("J9").Value = 0.4427387186280324
$ws.Range("O9").Value = 0.133770027207319
$ws.Range("P9").Value = 0.1849808797181794
$ws.Range("Q9").Value = 12.20835046950067
$ws.Range("R9").Value = 109.875154225506
$ws.Range("S9").Value = 0.05922517043660545
$ws.Range("T9").Value = 0.08189819765711295

# Row 10
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.2919476666666667
$ws.Range("H10").Value = 0.875843
$ws.Range("I10").Value = 0.4427387186280323
$ws.Range("J10").Value = 0.4427387186280324
$ws.Range("M10").Value = 3.307267666666667
$ws.Range("N10").Value = 9.921803
$ws.Range("O10").Value = 0.01057976888853842
$ws.Range("P10").Value = 0.01462999594956937
$ws.Range("Q10").Value = 0.9655490783254445
$ws.Range("R10").Value = 8.689941704929002
$ws.Range("S10").Value = 0.004684073321092219
$ws.Range("T10").Value = 0.006477265660245647

# Row 11
$ws.Range("D11").Value = "MuSCs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.2919476666666667
$ws.Range("H11").Value = 0.875843
$ws.Range("I11").Value = 0.4427387186280323
$ws.Range("J11").Value = 0.4427387186280324
$ws.Range("M11").Value = 18.767532
$ws.Range("N11").Value = 37.535064
$ws.Range("O11").Value = 0.06003631129389966
$ws.Range("P11").Value = 0.05534657705729765
$ws.Range("Q11").Value = 5.479137176491999
$ws.Range("R11").Value = 32.874823058952
$ws.Range("S11").Value = 0.0265803995334148
$ws.Range("T11").Value = 0.02450407260679562

# Row 12
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.2919476666666667
$ws.Range("H12").Value = 0.875843
$ws.Range("I12").Value = 0.4427387186280323
$ws.Range("J12").Value = 0.4427387186280324
$ws.Range("M12").Value = 1.967337333333333
$ws.Range("N12").Value = 5.902012
$ws.Range("O12").Value = 0.006293404831498911
$ws.Range("P12").Value = 0.008702693618721296
$ws.Range("Q12").Value = 0.5743595440128889
$ws.Range("R12").Value = 5.169235896116001
$ws.Range("S12").Value = 0.002786333990905295
$ws.Range("T12").Value = 0.003853019421365021

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.2919476666666667
$ws.Range("H13").Value = 0.875843
$ws.Range("I13").Value = 0.4427387186280323
$ws.Range("J13").Value = 0.4427387186280324
$ws.Range("M13").Value = 5.884659333333333
$ws.Range("N13").Value = 17.653978
$ws.Range("O13").Value = 0.01882470426023795
$ws.Range("P13").Value = 0.02603131977462027
$ws.Range("Q13").Value = 1.718012561494889
$ws.Range("R13").Value = 15.462113053454
$ws.Range("S13").Value = 0.008334425442729408
$ws.Range("T13").Value = 0.01152507316121194
